# Updates the cryptos list on the active sheet: refreshed Price (col D) and
# Volume(1h) (col E) figures, plus two row swaps caused by re-ranking
# (Chainlink/WrappedEther at rows 16-17, TrustWalletToken/Stacks at rows 50-51).
# Price values are prefixed with a leading apostrophe so Excel keeps them as
# literal text (matching the original formatting, e.g. trailing zeros like
# "87.20") instead of silently re-parsing them as numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '''39.949.28'
$ws.Range('E2').Value = '  -0.63%  '
$ws.Range('D3').Value = '''2.216.12'
$ws.Range('E3').Value = '  -0.91%  '
$ws.Range('E4').Value = '  -0.12%  '
$ws.Range('D5').Value = '''291.98'
$ws.Range('E5').Value = '  -0.41%  '
$ws.Range('D6').Value = '''87.20'
$ws.Range('E6').Value = '  -2.15%  '
$ws.Range('D7').Value = '''0.510'
$ws.Range('E7').Value = '  -1.12%  '
$ws.Range('E8').Value = '  +0.00%  '
$ws.Range('D9').Value = '''0.465'
$ws.Range('E9').Value = '  -2.02%  '
$ws.Range('D10').Value = '''30.44'
$ws.Range('E10').Value = '  -2.82%  '
$ws.Range('D11').Value = '''0.0778'
$ws.Range('E11').Value = '  -1.66%  '
$ws.Range('D12').Value = '''50.02'
$ws.Range('E12').Value = '  +5.05%  '
$ws.Range('D13').Value = '''0.113'
$ws.Range('E13').Value = '  +3.22%  '
$ws.Range('D14').Value = '''6.40'
$ws.Range('E14').Value = '  -0.93%  '
$ws.Range('D15').Value = '''2.556.91'
$ws.Range('E15').Value = '  -0.95%  '
$ws.Range('B16').Value = 'Chainlink'
$ws.Range('C16').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D16').Value = '''13.71'
$ws.Range('E16').Value = '  -3.53%  '
$ws.Range('B17').Value = 'WrappedEther'
$ws.Range('C17').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D17').Value = '''2.243.99'
$ws.Range('E17').Value = '  -0.43%  '
$ws.Range('D18').Value = '''0.730'
$ws.Range('E18').Value = '  -1.04%  '
$ws.Range('D19').Value = '''39.837.46'
$ws.Range('E19').Value = '  -0.77%  '
$ws.Range('D20').Value = '''0.0₃0884'
$ws.Range('E20').Value = '  -0.82%  '
$ws.Range('D21').Value = '''11.09'
$ws.Range('E21').Value = '  -6.70%  '
$ws.Range('D22').Value = '''5.73'
$ws.Range('E22').Value = '  -2.20%  '
$ws.Range('D23').Value = '''65.41'
$ws.Range('E23').Value = '  -0.91%  '
$ws.Range('D24').Value = '''235.84'
$ws.Range('E24').Value = '  +0.01%  '
$ws.Range('E25').Value = '  -0.02%  '
$ws.Range('D26').Value = '''2.45'
$ws.Range('E26').Value = '  -1.55%  '
$ws.Range('D27').Value = '''1.81'
$ws.Range('E27').Value = '  -4.17%  '
$ws.Range('D28').Value = '''2.33'
$ws.Range('E28').Value = '  +5.75%  '
$ws.Range('D29').Value = '''23.11'
$ws.Range('E29').Value = '  +0.52%  '
$ws.Range('D30').Value = '''9.21'
$ws.Range('E30').Value = '  -1.32%  '
$ws.Range('D31').Value = '''157.85'
$ws.Range('E31').Value = '  +3.51%  '
$ws.Range('D32').Value = '''31.47'
$ws.Range('E32').Value = '  -4.82%  '
$ws.Range('E33').Value = '  +0.02%  '
$ws.Range('D34').Value = '''4.94'
$ws.Range('E34').Value = '  -1.68%  '
$ws.Range('D35').Value = '''3.03'
$ws.Range('E35').Value = '  +5.24%  '
$ws.Range('D36').Value = '''0.0708'
$ws.Range('E36').Value = '  -2.26%  '
$ws.Range('E37').Value = '  -1.69%  '
$ws.Range('D38').Value = '''0.113'
$ws.Range('E38').Value = '  +0.12%  '
$ws.Range('D39').Value = '''0.0982'
$ws.Range('E39').Value = '  -1.98%  '
$ws.Range('D40').Value = '''1.72'
$ws.Range('E40').Value = '  -0.56%  '
$ws.Range('D41').Value = '''15.25'
$ws.Range('E41').Value = '  -6.16%  '
$ws.Range('D42').Value = '''2.093.17'
$ws.Range('E42').Value = '  -0.69%  '
$ws.Range('D43').Value = '''3.70'
$ws.Range('E43').Value = '  -5.08%  '
$ws.Range('D44').Value = '''0.0269'
$ws.Range('E44').Value = '  -0.98%  '
$ws.Range('D45').Value = '''17.86'
$ws.Range('E45').Value = '  -1.22%  '
$ws.Range('D46').Value = '''9.90'
$ws.Range('E46').Value = '  -1.67%  '
$ws.Range('D47').Value = '''1.99'
$ws.Range('E47').Value = '  -9.02%  '
$ws.Range('D48').Value = '''2.68'
$ws.Range('E48').Value = '  -1.46%  '
$ws.Range('D49').Value = '''2.428.43'
$ws.Range('E49').Value = '  -0.94%  '
$ws.Range('B50').Value = 'TrustWalletToken'
$ws.Range('C50').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D50').Value = '''1.11'
$ws.Range('E50').Value = '  +2.04%  '
$ws.Range('B51').Value = 'Stacks'
$ws.Range('C51').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D51').Value = '''1.45'
$ws.Range('E51').Value = '  -1.85%  '
